$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.599.07"
Set-TextValue "E2" "  -0.23%  "
Set-TextValue "D3" "1.753.62"
Set-TextValue "E3" "  -0.38%  "
Set-TextValue "D4" "1.004"
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "324.02"
Set-TextValue "E5" "  +0.03%  "
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.13%  "
Set-TextValue "E7" "  +4.78%  "
Set-TextValue "D8" "0.3549"
Set-TextValue "E8" "  -1.61%  "
Set-TextValue "D9" "0.07458"
Set-TextValue "E9" "  -1.38%  "
Set-TextValue "D10" "41.58"
Set-TextValue "E10" "  -1.33%  "
Set-TextValue "D11" "1.082"
Set-TextValue "E11" "  -2.59%  "
Set-TextValue "D12" "1.002"
Set-TextValue "E12" "  +0.06%  "
Set-TextValue "D13" "20.71"
Set-TextValue "E13" "  -0.35%  "
Set-TextValue "D14" "5.971"
Set-TextValue "E14" "  -1.73%  "
Set-TextValue "D15" "7.142"
Set-TextValue "E15" "  -1.44%  "
Set-TextValue "D16" "1.753.75"
Set-TextValue "E16" "  -0.36%  "
Set-TextValue "D17" "93.28"
Set-TextValue "E17" "  +0.85%  "
Set-TextValue "D18" "0.00001057"
Set-TextValue "E18" "  -0.99%  "
Set-TextValue "D19" "0.06471"
Set-TextValue "E19" "  +0.60%  "
Set-TextValue "E20" "  +0.07%  "
Set-TextValue "D21" "17.03"
Set-TextValue "D22" "5.748"
Set-TextValue "E22" "  -2.30%  "
Set-TextValue "D23" "27.664.76"
Set-TextValue "E23" "  -0.17%  "
Set-TextValue "D24" "11.21"
Set-TextValue "E24" "  -0.54%  "
Set-TextValue "E25" "  +0.74%  "
Set-TextValue "D26" "164.87"
Set-TextValue "E26" "  +1.26%  "
Set-TextValue "D27" "20.15"
Set-TextValue "E27" "  -2.13%  "
Set-TextValue "D28" "1.953.39"
Set-TextValue "E28" "  -0.37%  "
Set-TextValue "D29" "2.077"
Set-TextValue "E29" "  -3.62%  "
Set-TextValue "D30" "125.04"
Set-TextValue "E30" "  -0.66%  "
Set-TextValue "D31" "1.090"
Set-TextValue "E31" "  -0.79%  "
Set-TextValue "D33" "3.655"
Set-TextValue "E33" "  -0.67%  "
Set-TextValue "D34" "5.479"
Set-TextValue "E34" "  -2.17%  "
Set-TextValue "D35" "0.02285"
Set-TextValue "E35" "  -0.63%  "
Set-TextValue "D36" "11.67"
Set-TextValue "E36" "  -4.35%  "
Set-TextValue "D37" "0.06029"
Set-TextValue "E37" "  +0.30%  "
Set-TextValue "D38" "0.2079"
Set-TextValue "E38" "  -1.75%  "
Set-TextValue "D39" "0.6290"
Set-TextValue "E39" "  -1.04%  "
Set-TextValue "D40" "4.931"
Set-TextValue "E40" "  -0.56%  "
Set-TextValue "E41" "  -0.56%  "
Set-TextValue "D42" "1.391"
Set-TextValue "E42" "  -0.38%  "
Set-TextValue "D43" "7.741"
Set-TextValue "E43" "  -2.07%  "
Set-TextValue "D44" "13.17"
Set-TextValue "E44" "  -1.94%  "
Set-TextValue "D45" "3.714"
Set-TextValue "E45" "  -0.02%  "
Set-TextValue "D46" "0.5855"
Set-TextValue "E46" "  -1.33%  "
Set-TextValue "D47" "122.58"
Set-TextValue "E47" "  -0.71%  "
Set-TextValue "D48" "1.935"
Set-TextValue "E48" "  -2.85%  "
Set-TextValue "D49" "0.06893"
Set-TextValue "E49" "  +0.27%  "
Set-TextValue "D50" "1.126"
Set-TextValue "E50" "  -4.15%  "
Set-TextValue "D51" "71.61"
Set-TextValue "E51" "  -2.30%  "
